# Continue the schedule: replace the last block of 5-minute time slots
# (rows 7-12) with the next chronological block, and drop the now-removed
# trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "22:30-22:35"
$ws.Range("C8").Value = "22:35-22:40"
$ws.Range("C9").Value = "22:40-22:45"
$ws.Range("C10").Value = "22:45-22:50"
$ws.Range("C11").Value = "22:50-22:55"

# Row 12 (old "22:25-22:30" slot) is no longer part of the schedule.
$ws.Rows("12:12").Delete()

# Cursor/selection ends up on column B instead of C.
$ws.Range("B15").Select()
